# account.xlsx: add a second "username" column next to the existing
# password/gizem_guvel pair (new header "password2" + new value "gzm_gvl"),
# and drop the frozen header pane that no longer makes sense once the grid
# is being actively edited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unfreeze the header/column split that was previously pinned at B3.
$win = $excel.ActiveWindow
$win.FreezePanes = $false

# New column of data: header cell (row 2) + value cell (row 3).
$ws.Range("B2").Value = "password2"
$ws.Range("B3").Value = "gzm_gvl"

# Pick up the same look-and-feel as the existing A2/A3 cells so the new
# column reads as part of the same header/value rows instead of a stray
# unformatted cell.
$ws.Range("A2").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("A3").Copy()
$ws.Range("B3").PasteSpecial(-4122)
